# Update column C ("Fitness") values on Sheet1 to reflect the latest run's
# logged fitness numbers. The values plateau in several blocks as the run
# progresses through generations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (startRow, endRow, newFitnessValue) using 1-based sheet rows.
# Row 1 is the header (Run / Generation / Fitness); data starts at row 2.
$blocks = @(
    @{ Start = 2;   End = 17;  Value = 8129 },
    @{ Start = 18;  End = 24;  Value = 7912 },
    @{ Start = 25;  End = 25;  Value = 7836 },
    @{ Start = 26;  End = 212; Value = 7590 },
    @{ Start = 213; End = 252; Value = 7573 }
)

foreach ($block in $blocks) {
    for ($r = $block.Start; $r -le $block.End; $r++) {
        $ws.Cells.Item($r, 3).Value = $block.Value
    }
}
